# Update the Jogos_do_Dia_Betfair_Back_Lay sheet:
#  - Row 2 odds (F2:AO2) reset to 0
#  - Row 3: AO3 36 -> 38
#  - Row 4: league/time/teams updated (Welsh Premiership/New Saints-Bala Town
#    -> English Premier League/West Ham-Nottm Forest), odds (F4:AO4) reset to 0
#  - Row 5: league/teams updated (English Premier League/West Ham-Nottm Forest
#    -> Scottish Premiership/Rangers-Aberdeen), odds (F5:AO5) reset to 0
#  - Row 6 (old Scottish Premiership/Rangers-Aberdeen row) deleted entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: zero out odds columns F:AO
$ws.Range("F2:AO2").Value = 0

# Row 3: only AO3 changes
$ws.Range("AO3").Value = 38

# Row 4: new match data (was Welsh Premiership / The New Saints vs Bala Town)
$ws.Range("A4").Value = "English Premier League"
$ws.Range("C4").Value = "17:00:00"
$ws.Range("D4").Value = "West Ham"
$ws.Range("E4").Value = "Nottm Forest"
$ws.Range("F4:AO4").Value = 0

# Row 5: new match data (was English Premier League / West Ham vs Nottm Forest)
$ws.Range("A5").Value = "Scottish Premiership"
$ws.Range("D5").Value = "Rangers"
$ws.Range("E5").Value = "Aberdeen"
$ws.Range("F5:AO5").Value = 0

# Row 6: remove entirely (old Scottish Premiership / Rangers vs Aberdeen row)
$ws.Rows("6").Delete()
